# Insert a new data row at row 74 (pushing the existing rows 74-198 down to 75-199)
# and populate it with the new record, matching the values in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 74; Excel shifts rows 74..198 down to 75..199.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new record's data.
$ws.Cells.Item(74, 1).Value = 4
$ws.Cells.Item(74, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(74, 3).Value = "Los Lagos"
$ws.Cells.Item(74, 4).Value = 44533
$ws.Cells.Item(74, 5).Value = 10
$ws.Cells.Item(74, 6).Value = 100112037
$ws.Cells.Item(74, 7).Value = "Cebollín"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 140
$ws.Cells.Item(74, 11).Value = 6000
$ws.Cells.Item(74, 12).Value = 6000
$ws.Cells.Item(74, 13).Value = 6000
$ws.Cells.Item(74, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(74, 15).Value = "Región Metropolitana"
$ws.Cells.Item(74, 16).Value = 167
$ws.Cells.Item(74, 17).Value = 36
$ws.Cells.Item(74, 18).Value = "Hortaliza"

# Match the date cell's style/number format to the rest of column D.
$ws.Cells.Item(74, 4).NumberFormat = $ws.Cells.Item(75, 4).NumberFormat
